# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Column D ("Price") and column E ("Volume(1h)") are plain text cells (inlineStr)
# in the source file, so number-looking prices are written with a leading
# apostrophe to keep Excel from auto-coercing them to numeric cells, and the
# style is reset to "Normal" right after so no stray text-format style sticks
# around on the cell (matches the original, unstyled text cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.347.98'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '2.603.73'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''592.45'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.53%  '
$ws.Range('D6').Value = '''150.34'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.71%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -1.06%  '
$ws.Range('D9').Value = '2.603.15'
$ws.Range('E9').Value = '  -0.38%  '
$ws.Range('D10').Value = '''0.129'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2.16%  '
$ws.Range('E11').Value = '  -0.04%  '
$ws.Range('E12').Value = '  -1.43%  '
$ws.Range('E13').Value = '  -2.97%  '
$ws.Range('E14').Value = '  -2.90%  '
$ws.Range('D15').Value = '3.077.10'
$ws.Range('E15').Value = '  -0.38%  '
$ws.Range('E16').Value = '  -2.60%  '
$ws.Range('D17').Value = '67.192.43'
$ws.Range('E17').Value = '  -0.55%  '
$ws.Range('D18').Value = '2.605.73'
$ws.Range('E18').Value = '  -0.22%  '
$ws.Range('D19').Value = '''370.31'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.15%  '
$ws.Range('E20').Value = '  -2.16%  '
$ws.Range('D21').Value = '''7.34'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.23%  '
$ws.Range('E22').Value = '  -2.83%  '
$ws.Range('D23').Value = '''4.76'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -4.25%  '
$ws.Range('D24').Value = '''2.02'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -3.29%  '
$ws.Range('D25').Value = '''73.16'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +4.49%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  -1.58%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').Value = '''577.24'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.93%  '
$ws.Range('E31').Value = '  -6.00%  '
$ws.Range('E32').Value = '  -5.23%  '
$ws.Range('D33').Value = '''7.65'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -3.55%  '
$ws.Range('E34').Value = '  -2.99%  '
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('E36').Value = '  -4.19%  '
$ws.Range('D37').Value = '''1.50'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.44%  '
$ws.Range('D38').Value = '''158.39'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.84%  '
$ws.Range('D39').Value = '''19.03'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.92%  '
$ws.Range('E40').Value = '  +0.81%  '
$ws.Range('E41').Value = '  -1.87%  '
$ws.Range('E42').Value = '  -3.44%  '
$ws.Range('E43').Value = '  +4.13%  '
$ws.Range('E44').Value = '  -4.39%  '
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').Value = '''152.81'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.37%  '
$ws.Range('E47').Value = '  -1.70%  '
$ws.Range('D48').Value = '''3.64'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.23%  '
$ws.Range('D49').Value = '''0.0776'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.56%  '
$ws.Range('E50').Value = '  -4.55%  '
$ws.Range('D51').Value = '''21.28'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.34%  '
